$d = $word.ActiveDocument

# The only real textual change in this edit is a typo fix: "wisse dood"
# becomes "gewisse dood" (i.e. a "ge" prefix is inserted before "wisse"),
# in the paragraph "Hij had hem willen overtuigen... op de berg te
# wachten op een wisse dood."
$full = $d.Content.Text
$needle = "op de berg te wachten op een wisse dood"
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $prefixLen = ("op de berg te wachten op een ").Length
    $pos = $idx + $prefixLen
    $insertPoint = $d.Range($pos, $pos)
    $insertPoint.InsertBefore("ge")
}
